# Apply NATMI edit: add "ECs" as a third sending cluster (rows 2-16)
# Updates existing rows 2-11 and appends new rows 12-16 for S100a9 -> Tlr4 (ECs sender)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "M1"
$ws.Range("B2").Value = "S100a9"
$ws.Range("C2").Value = "Tlr4"
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 1
$ws.Range("F2").Value = 0.3333333333333333
$ws.Range("G2").Value = 0.08871633333333334
$ws.Range("H2").Value = 0.266149
$ws.Range("I2").Value = 0.0007521739725275529
$ws.Range("J2").Value = 0.0007521739725275529
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 9.115710333333334
$ws.Range("N2").Value = 27.347131
$ws.Range("O2").Value = 0.1106357638930427
$ws.Range("P2").Value = 0.1106357638930427
$ws.Range("Q2").Value = 0.8087123965021112
$ws.Range("R2").Value = 7.278411568519001
$ws.Range("S2").Value = 0.00008321734203105032
$ws.Range("T2").Value = 0.0000832173420310503

$ws.Range("A3").Value = "M1"
$ws.Range("B3").Value = "S100a9"
$ws.Range("C3").Value = "Tlr4"
$ws.Range("D3").Value = "FAPs"
$ws.Range("E3").Value = 1
$ws.Range("F3").Value = 0.3333333333333333
$ws.Range("G3").Value = 0.08871633333333334
$ws.Range("H3").Value = 0.266149
$ws.Range("I3").Value = 0.0007521739725275529
$ws.Range("J3").Value = 0.0007521739725275529
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 9.004096
$ws.Range("N3").Value = 27.012288
$ws.Range("O3").Value = 0.1092811204721574
$ws.Range("P3").Value = 0.1092811204721574
$ws.Range("Q3").Value = 0.7988103821013335
$ws.Range("R3").Value = 7.189293438912001
$ws.Range("S3").Value = 0.00008219841450780471
$ws.Range("T3").Value = 0.00008219841450780471

$ws.Range("A4").Value = "M1"
$ws.Range("B4").Value = "S100a9"
$ws.Range("C4").Value = "Tlr4"
$ws.Range("D4").Value = "M1"
$ws.Range("E4").Value = 1
$ws.Range("F4").Value = 0.3333333333333333
$ws.Range("G4").Value = 0.08871633333333334
$ws.Range("H4").Value = 0.266149
$ws.Range("I4").Value = 0.0007521739725275529
$ws.Range("J4").Value = 0.0007521739725275529
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 31.76332533333333
$ws.Range("N4").Value = 95.289976
$ws.Range("O4").Value = 0.3855058611490069
$ws.Range("P4").Value = 0.3855058611490069
$ws.Range("Q4").Value = 2.817925758047111
$ws.Range("R4").Value = 25.361331822424
$ws.Range("S4").Value = 0.0002899674750131037
$ws.Range("T4").Value = 0.0002899674750131037

$ws.Range("A5").Value = "M1"
$ws.Range("B5").Value = "S100a9"
$ws.Range("C5").Value = "Tlr4"
$ws.Range("D5").Value = "M2"
$ws.Range("E5").Value = 1
$ws.Range("F5").Value = 0.3333333333333333
$ws.Range("G5").Value = 0.08871633333333334
$ws.Range("H5").Value = 0.266149
$ws.Range("I5").Value = 0.0007521739725275529
$ws.Range("J5").Value = 0.0007521739725275529
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 30.12843966666667
$ws.Range("N5").Value = 90.385319
$ws.Range("O5").Value = 0.3656635429976674
$ws.Range("P5").Value = 0.3656635429976673
$ws.Range("Q5").Value = 2.672884696281222
$ws.Range("R5").Value = 24.055962266531
$ws.Range("S5").Value = 0.0002750425997450551
$ws.Range("T5").Value = 0.0002750425997450551

$ws.Range("A6").Value = "M1"
$ws.Range("B6").Value = "S100a9"
$ws.Range("C6").Value = "Tlr4"
$ws.Range("D6").Value = "sCs"
$ws.Range("E6").Value = 1
$ws.Range("F6").Value = 0.3333333333333333
$ws.Range("G6").Value = 0.08871633333333334
$ws.Range("H6").Value = 0.266149
$ws.Range("I6").Value = 0.0007521739725275529
$ws.Range("J6").Value = 0.0007521739725275529
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 2.382313
$ws.Range("N6").Value = 7.146939
$ws.Range("O6").Value = 0.0289137114881257
$ws.Range("P6").Value = 0.02891371148812569
$ws.Range("Q6").Value = 0.2113500742123333
$ws.Range("R6").Value = 1.902150667911
$ws.Range("S6").Value = 0.00002174814123053905
$ws.Range("T6").Value = 0.00002174814123053905

$ws.Range("A7").Value = "M2"
$ws.Range("B7").Value = "S100a9"
$ws.Range("C7").Value = "Tlr4"
$ws.Range("D7").Value = "ECs"
$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 41.455644
$ws.Range("H7").Value = 124.366932
$ws.Range("I7").Value = 0.3514781918906479
$ws.Range("J7").Value = 0.3514781918906478
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 9.115710333333334
$ws.Range("N7").Value = 27.347131
$ws.Range("O7").Value = 0.1106357638930427
$ws.Range("P7").Value = 0.1106357638930427
$ws.Range("Q7").Value = 377.8976423857881
$ws.Range("R7").Value = 3401.078781472092
$ws.Range("S7").Value = 0.03888605825156727
$ws.Range("T7").Value = 0.03888605825156725

$ws.Range("A8").Value = "M2"
$ws.Range("B8").Value = "S100a9"
$ws.Range("C8").Value = "Tlr4"
$ws.Range("D8").Value = "FAPs"
$ws.Range("E8").Value = 3
$ws.Range("F8").Value = 1
$ws.Range("G8").Value = 41.455644
$ws.Range("H8").Value = 124.366932
$ws.Range("I8").Value = 0.3514781918906479
$ws.Range("J8").Value = 0.3514781918906478
$ws.Range("K8").Value = 3
$ws.Range("L8").Value = 1
$ws.Range("M8").Value = 9.004096
$ws.Range("N8").Value = 27.012288
$ws.Range("O8").Value = 0.1092811204721574
$ws.Range("P8").Value = 0.1092811204721574
$ws.Range("Q8").Value = 373.270598317824
$ws.Range("R8").Value = 3359.435384860416
$ws.Range("S8").Value = 0.03840993063133794
$ws.Range("T8").Value = 0.03840993063133794

$ws.Range("A9").Value = "M2"
$ws.Range("B9").Value = "S100a9"
$ws.Range("C9").Value = "Tlr4"
$ws.Range("D9").Value = "M1"
$ws.Range("E9").Value = 3
$ws.Range("F9").Value = 1
$ws.Range("G9").Value = 41.455644
$ws.Range("H9").Value = 124.366932
$ws.Range("I9").Value = 0.3514781918906479
$ws.Range("J9").Value = 0.3514781918906478
$ws.Range("K9").Value = 3
$ws.Range("L9").Value = 1
$ws.Range("M9").Value = 31.76332533333333
$ws.Range("N9").Value = 95.289976
$ws.Range("O9").Value = 0.3855058611490069
$ws.Range("P9").Value = 0.3855058611490069
$ws.Range("Q9").Value = 1316.769107274848
$ws.Range("R9").Value = 11850.92196547363
$ws.Range("S9").Value = 0.1354969030399001
$ws.Range("T9").Value = 0.1354969030399001

$ws.Range("A10").Value = "M2"
$ws.Range("B10").Value = "S100a9"
$ws.Range("C10").Value = "Tlr4"
$ws.Range("D10").Value = "M2"
$ws.Range("E10").Value = 3
$ws.Range("F10").Value = 1
$ws.Range("G10").Value = 41.455644
$ws.Range("H10").Value = 124.366932
$ws.Range("I10").Value = 0.3514781918906479
$ws.Range("J10").Value = 0.3514781918906478
$ws.Range("K10").Value = 3
$ws.Range("L10").Value = 1
$ws.Range("M10").Value = 30.12843966666667
$ws.Range("N10").Value = 90.385319
$ws.Range("O10").Value = 0.3656635429976674
$ws.Range("P10").Value = 0.3656635429976673
$ws.Range("Q10").Value = 1248.993869096812
$ws.Range("R10").Value = 11240.94482187131
$ws.Range("S10").Value = 0.1285227609331483
$ws.Range("T10").Value = 0.1285227609331483

$ws.Range("A11").Value = "M2"
$ws.Range("B11").Value = "S100a9"
$ws.Range("C11").Value = "Tlr4"
$ws.Range("D11").Value = "sCs"
$ws.Range("E11").Value = 3
$ws.Range("F11").Value = 1
$ws.Range("G11").Value = 41.455644
$ws.Range("H11").Value = 124.366932
$ws.Range("I11").Value = 0.3514781918906479
$ws.Range("J11").Value = 0.3514781918906478
$ws.Range("K11").Value = 3
$ws.Range("L11").Value = 1
$ws.Range("M11").Value = 2.382313
$ws.Range("N11").Value = 7.146939
$ws.Range("O11").Value = 0.0289137114881257
$ws.Range("P11").Value = 0.02891371148812569
$ws.Range("Q11").Value = 98.76031962457199
$ws.Range("R11").Value = 888.8428766211479
$ws.Range("S11").Value = 0.01016253903469427
$ws.Range("T11").Value = 0.01016253903469427

$ws.Range("A12").Value = "ECs"
$ws.Range("B12").Value = "S100a9"
$ws.Range("C12").Value = "Tlr4"
$ws.Range("D12").Value = "ECs"
$ws.Range("E12").Value = 3
$ws.Range("F12").Value = 1
$ws.Range("G12").Value = 76.40220066666666
$ws.Range("H12").Value = 229.206602
$ws.Range("I12").Value = 0.6477696341368246
$ws.Range("J12").Value = 0.6477696341368245
$ws.Range("K12").Value = 3
$ws.Range("L12").Value = 1
$ws.Range("M12").Value = 9.115710333333334
$ws.Range("N12").Value = 27.347131
$ws.Range("O12").Value = 0.1106357638930427
$ws.Range("P12").Value = 0.1106357638930427
$ws.Range("Q12").Value = 696.4603301065403
$ws.Range("R12").Value = 6268.142970958862
$ws.Range("S12").Value = 0.07166648829944437
$ws.Range("T12").Value = 0.07166648829944434

$ws.Range("A13").Value = "ECs"
$ws.Range("B13").Value = "S100a9"
$ws.Range("C13").Value = "Tlr4"
$ws.Range("D13").Value = "FAPs"
$ws.Range("E13").Value = 3
$ws.Range("F13").Value = 1
$ws.Range("G13").Value = 76.40220066666666
$ws.Range("H13").Value = 229.206602
$ws.Range("I13").Value = 0.6477696341368246
$ws.Range("J13").Value = 0.6477696341368245
$ws.Range("K13").Value = 3
$ws.Range("L13").Value = 1
$ws.Range("M13").Value = 9.004096
$ws.Range("N13").Value = 27.012288
$ws.Range("O13").Value = 0.1092811204721574
$ws.Range("P13").Value = 0.1092811204721574
$ws.Range("Q13").Value = 687.9327494139307
$ws.Range("R13").Value = 6191.394744725376
$ws.Range("S13").Value = 0.07078899142631165
$ws.Range("T13").Value = 0.07078899142631163

$ws.Range("A14").Value = "ECs"
$ws.Range("B14").Value = "S100a9"
$ws.Range("C14").Value = "Tlr4"
$ws.Range("D14").Value = "M1"
$ws.Range("E14").Value = 3
$ws.Range("F14").Value = 1
$ws.Range("G14").Value = 76.40220066666666
$ws.Range("H14").Value = 229.206602
$ws.Range("I14").Value = 0.6477696341368246
$ws.Range("J14").Value = 0.6477696341368245
$ws.Range("K14").Value = 3
$ws.Range("L14").Value = 1
$ws.Range("M14").Value = 31.76332533333333
$ws.Range("N14").Value = 95.289976
$ws.Range("O14").Value = 0.3855058611490069
$ws.Range("P14").Value = 0.3855058611490069
$ws.Range("Q14").Value = 2426.78795595795
$ws.Range("R14").Value = 21841.09160362155
$ws.Range("S14").Value = 0.2497189906340937
$ws.Range("T14").Value = 0.2497189906340937

$ws.Range("A15").Value = "ECs"
$ws.Range("B15").Value = "S100a9"
$ws.Range("C15").Value = "Tlr4"
$ws.Range("D15").Value = "M2"
$ws.Range("E15").Value = 3
$ws.Range("F15").Value = 1
$ws.Range("G15").Value = 76.40220066666666
$ws.Range("H15").Value = 229.206602
$ws.Range("I15").Value = 0.6477696341368246
$ws.Range("J15").Value = 0.6477696341368245
$ws.Range("K15").Value = 3
$ws.Range("L15").Value = 1
$ws.Range("M15").Value = 30.12843966666667
$ws.Range("N15").Value = 90.385319
$ws.Range("O15").Value = 0.3656635429976674
$ws.Range("P15").Value = 0.3656635429976673
$ws.Range("Q15").Value = 2301.879093186226
$ws.Range("R15").Value = 20716.91183867603
$ws.Range("S15").Value = 0.2368657394647741
$ws.Range("T15").Value = 0.236865739464774

$ws.Range("A16").Value = "ECs"
$ws.Range("B16").Value = "S100a9"
$ws.Range("C16").Value = "Tlr4"
$ws.Range("D16").Value = "sCs"
$ws.Range("E16").Value = 3
$ws.Range("F16").Value = 1
$ws.Range("G16").Value = 76.40220066666666
$ws.Range("H16").Value = 229.206602
$ws.Range("I16").Value = 0.6477696341368246
$ws.Range("J16").Value = 0.6477696341368245
$ws.Range("K16").Value = 3
$ws.Range("L16").Value = 1
$ws.Range("M16").Value = 2.382313
$ws.Range("N16").Value = 7.146939
$ws.Range("O16").Value = 0.0289137114881257
$ws.Range("P16").Value = 0.02891371148812569
$ws.Range("Q16").Value = 182.0139558768086
$ws.Range("R16").Value = 1638.125602891278
$ws.Range("S16").Value = 0.01872942431220089
$ws.Range("T16").Value = 0.01872942431220088

